$wb = $excel.ActiveWorkbook

# Sheet 1 = 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value  = 13381
$ws1.Range("F5").Value  = 772
$ws1.Range("F13").Value = 21136
$ws1.Range("G13").Value = 0
$ws1.Range("F14").Value = 544
$ws1.Range("F16").Value = 513
$ws1.Range("F23").Value = 34
$ws1.Range("F25").Value = 284
$ws1.Range("F28").Value = 64
$ws1.Range("F29").Value = 380

# Sheet 2 = 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value  = 301
$ws2.Range("F4").Value  = 4479
$ws2.Range("F11").Value = 390

# Sheet 3 = 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 4437
$ws3.Range("F4").Value = 106

# Sheet 4 = 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F6").Value  = 13381
$ws4.Range("F7").Value  = 301
$ws4.Range("F8").Value  = 772
$ws4.Range("F9").Value  = 4437
$ws4.Range("F16").Value = 106
$ws4.Range("F17").Value = 21136
$ws4.Range("G17").Value = 0
$ws4.Range("F18").Value = 544
$ws4.Range("F19").Value = 4479
$ws4.Range("F23").Value = 513
$ws4.Range("F29").Value = 390
$ws4.Range("F35").Value = 34
$ws4.Range("F40").Value = 284
$ws4.Range("F43").Value = 64
$ws4.Range("F45").Value = 380
